$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add three new prevalence-output columns (all_cases, general_prob,
#    whole_pop) to the "cases2010district" sheet.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("cases2010district")

# Header row. Written in this order so the shared-string table gets the
# new unique strings in the same order as the target workbook:
# all_cases, general_prob, whole_pop.
$ws6.Range("F1").Value = "all_cases"
$ws6.Range("H1").Value = "general_prob"
$ws6.Range("G1").Value = "whole_pop"

# F2:F33 -- constant "all_cases" figure used for every district.
$ws6.Range("F2:F33").Value = 50000

# G2 -- total population across all districts (sum of pop column);
# G3:G33 -- same number, stored as a literal (matches the source
# workbook's pattern of repeating a computed constant down the column).
$ws6.Range("G2").Formula = "=SUM(D2:D33)"
$ws6.Range("G3:G33").Value = 17563749

# H2 -- standalone formula (mirrors the existing E2 pattern, which is also
# not part of the shared-formula group below it);
# H3:H33 -- shared formula group, like E3:E33.
$ws6.Range("H2").Formula = "=F2/G2"
$ws6.Range("H3:H33").Formula = "=F3/G3"

# ---------------------------------------------------------------------------
# 2. View changes: "cases2010district" becomes the selected/active tab, its
#    frozen pane moves from column A to column B, and the active cell
#    becomes the newly added H1 header.
# ---------------------------------------------------------------------------
$ws6.Activate() | Out-Null
$ws6.Range("B1").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
$ws6.Range("H1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. "WHO_estimates" loses tab-selected status; its frozen pane moves from
#    column L to column B and the selection moves to C12:E12.
# ---------------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("WHO_estimates")
$ws12.Activate() | Out-Null
$ws12.Range("L1").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
$ws12.Range("C12:E12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Re-activate "cases2010district" last so it ends up as the workbook's
#    active tab (matches workbookView activeTab = 5, the sheet's 0-based
#    position).
# ---------------------------------------------------------------------------
$ws6.Activate() | Out-Null
$ws6.Range("H1").Select() | Out-Null
